# Update data_source column (D2:D5) from "recovered" to "recovered_host"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "recovered_host"
$ws.Range("D3").Value = "recovered_host"
$ws.Range("D4").Value = "recovered_host"
$ws.Range("D5").Value = "recovered_host"

# Update the active selection, matching the new view state in the diff
$ws.Range("C20").Select()
